$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 697, shifting existing rows 697:738 down to 698:739
$ws.Rows.Item(697).Insert()

# Populate the newly inserted row with the new data
$ws.Cells.Item(697, 1).NumberFormat = "@"
$ws.Cells.Item(697, 1).Value = "2026/01/21"
$ws.Cells.Item(697, 1).ClearFormats()
$ws.Cells.Item(697, 2).Value = "水"
$ws.Cells.Item(697, 3).Value = 16
$ws.Cells.Item(697, 4).Value = 201
